$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.343.41"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "1.879.34"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7112"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.48"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08028"
$ws.Range("E8").Value = "  +3.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3164"
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08308"
$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("D12").Value = "1.880.32"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.254"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.50"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7155"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.354"
$ws.Range("E16").Value = "  +4.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008543"
$ws.Range("E17").Value = "  +3.77%  "

$ws.Range("D18").Value = "29.355.64"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.56"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").Value = "2.132.84"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.808"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1562"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.083"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.88"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.509"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.327"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  -7.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05385"
$ws.Range("E33").Value = "  +1.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.937"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7706"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.187"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").Value = "1.263.01"
$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.509"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.30"
$ws.Range("E42").Value = "  +2.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9059"
$ws.Range("E43").Value = "  +1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.33"
$ws.Range("E44").Value = "  +1.90%  "

$ws.Range("E45").Value = "  +7.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").Value = "2.025.76"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5228"
$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.801"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.465"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4380"
$ws.Range("E51").Value = "  +1.41%  "
